$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly rows pulled from the Strava export for the Kilimanjaro tracker.
$rows = @(
    @("Matt", 45500, "Run", 53, 4.78, 253, 1, 38, 7, 1, 0, "Agile Antelope", 7),
    @("Matt", 45500, "Walk", 2, 0.08, 0, 2, 0, 0, 0, 0, "Agile Antelope", 7),
    @("Jeremiah", 45500, "Ride", 144, 39.98, 2243, 1, 36, 74, 29, 0, "Sauntering Hippo", 7),
    @("Matt", 45500, "Ride", 20, 6.11, 0, 3, 17, 0, 0, 0, "Agile Antelope", 7),
    @("Phil", 45500, "Walk", 80, 2.88, 345, 80, 0, 0, 0, 0, "Sauntering Hippo", 7),
    @("Steven", 45500, "Walk", 45, 2.23, 75, 45, 0, 0, 0, 0, "Brave Leopard", 7),
    @("Steven", 45500, "Walk", 27, 1.15, 30, 27, 0, 0, 0, 0, "Brave Leopard", 7),
    @("Eric", 45500, "Workout", 60, 0, 0, 21, 41, 0, 0, 0, "Sauntering Hippo", 7),
    @("Steven", 45501, "Walk", 33, 1.71, 79, 33, 0, 0, 0, 0, "Brave Leopard", 7),
    @("Steven", 45501, "Workout", 24, 0, 0, 16, 8, 0, 0, 0, "Brave Leopard", 7),
    @("Matt", 45501, "Run", 113, 10.01, 486, 4, 73, 21, 4, 0, "Agile Antelope", 7),
    @("Matt", 45501, "Walk", 6, 0.19, 0, 6, 0, 0, 0, 0, "Agile Antelope", 7)
)

$startRow = 290

# Copy the date format (numFmt applied to col B) from the previous row down onto
# the new date cells so the new rows render as dates like the rest of the table.
$ws.Cells.Item($startRow - 1, 2).Copy()
$lastRow = $startRow + $rows.Count - 1
$fmtTarget = $ws.Range($ws.Cells.Item($startRow, 2), $ws.Cells.Item($lastRow, 2))
$fmtTarget.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$r = $startRow
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

[void]$ws.Cells.Item($lastRow + 1, 1).Select()
